$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "51.904.39"
Set-TextValue $ws.Range("E2") "  +0.09%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.842.74"
Set-TextValue $ws.Range("E3") "  +1.73%  "

# Row 4
Set-TextValue $ws.Range("E4") "  -0.04%  "

# Row 5
Set-TextValue $ws.Range("D5") "350.65"
Set-TextValue $ws.Range("E5") "  -1.15%  "

# Row 6
Set-TextValue $ws.Range("D6") "113.52"
Set-TextValue $ws.Range("E6") "  +3.53%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.558"
Set-TextValue $ws.Range("E7") "  +1.18%  "

# Row 8
Set-TextValue $ws.Range("D8") "1.00"
Set-TextValue $ws.Range("E8") "  -0.02%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.620"
Set-TextValue $ws.Range("E9") "  +3.75%  "

# Row 10
Set-TextValue $ws.Range("D10") "40.38"
Set-TextValue $ws.Range("E10") "  +0.81%  "

# Row 11
Set-TextValue $ws.Range("E11") "  -0.92%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.0848"

# Row 13
Set-TextValue $ws.Range("D13") "20.11"
Set-TextValue $ws.Range("E13") "  -0.59%  "

# Row 14
Set-TextValue $ws.Range("D14") "7.80"
Set-TextValue $ws.Range("E14") "  +1.88%  "

# Row 15
Set-TextValue $ws.Range("D15") "3.272.75"
Set-TextValue $ws.Range("E15") "  +1.20%  "

# Row 16
Set-TextValue $ws.Range("D16") "0.987"
Set-TextValue $ws.Range("E16") "  +6.11%  "

# Row 17
Set-TextValue $ws.Range("D17") "2.839.24"
Set-TextValue $ws.Range("E17") "  +0.38%  "

# Row 18
Set-TextValue $ws.Range("D18") "51.873.87"
Set-TextValue $ws.Range("E18") "  +0.10%  "

# Row 20
Set-TextValue $ws.Range("D20") "7.66"
Set-TextValue $ws.Range("E20") "  -1.01%  "

# Row 21
Set-TextValue $ws.Range("D21") "13.51"
Set-TextValue $ws.Range("E21") "  +2.39%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.0₃0974"
Set-TextValue $ws.Range("E22") "  +0.73%  "

# Row 23
Set-TextValue $ws.Range("D23") "70.49"
Set-TextValue $ws.Range("E23") "  +0.52%  "

# Row 24
Set-TextValue $ws.Range("D24") "269.20"
Set-TextValue $ws.Range("E24") "  +0.82%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.76"
Set-TextValue $ws.Range("E25") "  +0.70%  "

# Row 26
Set-TextValue $ws.Range("D26") "26.35"
Set-TextValue $ws.Range("E26") "  +0.44%  "

# Row 27
Set-TextValue $ws.Range("E27") "  -0.05%  "

# Row 28
Set-TextValue $ws.Range("E28") "  +0.89%  "

# Row 29
Set-TextValue $ws.Range("D29") "39.62"
Set-TextValue $ws.Range("E29") "  +6.78%  "

# Row 30
Set-TextValue $ws.Range("D30") "10.57"
Set-TextValue $ws.Range("E30") "  +3.01%  "

# Row 31
Set-TextValue $ws.Range("E31") "  +16.06%  "

# Row 32
Set-TextValue $ws.Range("E32") "  +1.36%  "

# Row 33
Set-TextValue $ws.Range("D33") "52.79"
Set-TextValue $ws.Range("E33") "  +1.50%  "

# Row 34
Set-TextValue $ws.Range("E34") "  +1.31%  "

# Row 35
Set-TextValue $ws.Range("D35") "0.0896"
Set-TextValue $ws.Range("E35") "  +7.75%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.0451"
Set-TextValue $ws.Range("E36") "  -1.11%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.999"
Set-TextValue $ws.Range("E37") "  -0.11%  "

# Row 38
Set-TextValue $ws.Range("D38") "18.98"
Set-TextValue $ws.Range("E38") "  +2.09%  "

# Row 39
Set-TextValue $ws.Range("D39") "3.23"
Set-TextValue $ws.Range("E39") "  +2.08%  "

# Row 40
Set-TextValue $ws.Range("D40") "2.02"
Set-TextValue $ws.Range("E40") "  +2.32%  "

# Row 41
Set-TextValue $ws.Range("E41") "  +1.19%  "

# Row 42
Set-TextValue $ws.Range("D42") "2.54"
Set-TextValue $ws.Range("E42") "  -1.11%  "

# Row 43
Set-TextValue $ws.Range("D43") "123.16"
Set-TextValue $ws.Range("E43") "  +1.62%  "

# Row 44
Set-TextValue $ws.Range("B44") "EnergySwap"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D44") "22.32"
Set-TextValue $ws.Range("E44") "  +0.78%  "

# Row 45
Set-TextValue $ws.Range("B45") "WEMIXToken"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D45") "2.23"
Set-TextValue $ws.Range("E45") "  +1.72%  "

# Row 46
Set-TextValue $ws.Range("D46") "3.53"
Set-TextValue $ws.Range("E46") "  +6.77%  "

# Row 47
Set-TextValue $ws.Range("D47") "2.53"
Set-TextValue $ws.Range("E47") "  +8.90%  "

# Row 48
Set-TextValue $ws.Range("D48") "2.177.20"
Set-TextValue $ws.Range("E48") "  +1.61%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.251"
Set-TextValue $ws.Range("E49") "  +22.54%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.953"
Set-TextValue $ws.Range("E50") "  +4.21%  "

# Row 51
Set-TextValue $ws.Range("E51") "  -0.64%  "
